$d = $word.ActiveDocument

# 1. Update the date.
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2. Split the address paragraph "929 Story Road, San Jose CA 95122" into
#    two paragraphs: "929 Story Road" and "San Jose, CA 95122" (each keeping
#    the original paragraph/run formatting).
$addrRng = $d.Content
$addrRng.Find.Execute("929 Story Road, San Jose CA 95122", $true, $false, $false, $false, $false,
                       $true, 1, $false, "", 0)
$addrRng.InsertParagraphAfter()
$addrRng.Text = "929 Story Road"

$cityRng = $d.Content
$cityRng.Find.Execute("929 Story Road", $true, $false, $false, $false, $false,
                       $true, 1, $false, "", 0)
$addrPara = $cityRng.Paragraphs(1)
$cityPara = $addrPara.Next()
$cityPara.Range.InsertBefore("San Jose, CA 95122")

# 3. Remove the empty "NoSpacing" paragraph directly after "Board of Directors".
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.EndsWith("Board of Directors`r")) {
        $next = $p.Next()
        if ($next -ne $null -and $next.Range.Text -eq "`r" -and
            $next.Range.ParagraphFormat.Style.NameLocal -eq "No Spacing") {
            $next.Range.Delete()
        }
        break
    }
}
